# Applies the "refined workflow logic, tracked apps are now configured from
# storage bucket" commit: adds two new settings rows to the Settings sheet
# describing the TrackedApps file/storage-bucket, and removes the now
# obsolete TrackedApps sheet entirely.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new settings rows on the "Settings" sheet -------------------
$settings = $wb.Worksheets.Item("Settings")

$settings.Range("A6").Value = "TrackedAppsFileName"
$settings.Range("B6").Value = "TrackedApps.xlsx"
$settings.Range("C6").Value = "Name of file in Orchestrator storage that contains settings for which apps to track."

$settings.Range("A7").Value = "TrackedAppsStorageBucket"
$settings.Range("B7").Value = "TrackedAppsSettings"
$settings.Range("C7").Value = "Name of storage bucket in Orchestrator which contains the TrackedApps file"

$settings.Range("B13").Select()

# --- 2. Remove the obsolete "TrackedApps" sheet ------------------------------
$trackedApps = $wb.Worksheets.Item("TrackedApps")
$trackedApps.Delete()

# --- 3. Make "Assets" the active sheet --------------------------------------
$constants = $wb.Worksheets.Item("Constants")
$constants.Activate()
$constants.Range("A17").Select()

$assets = $wb.Worksheets.Item("Assets")
$assets.Activate()
